$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D8").Value = -8.773599999999989
$ws.Range("B12").Value = 5.565299999999999
$ws.Range("D12").Value = -7.953999999999997
$ws.Range("D14").Value = -8.693899999999998
$ws.Range("D22").Value = -7.900699999999995
